# Convert the "m:'doc.html'.fromHTMLURI()" field (fldChar/instrText based)
# into plain literal-text runs that spell out the field code itself, i.e.
# "{m:'doc.html'.fromHTMLURI()}", using <w:t> runs instead of <w:instrText>,
# and drop the now-unneeded leading/trailing space runs and the fldChar
# begin/end delimiters. The bookmark around "doc.html" / "'.fromHTMLURI()"
# is left untouched in place.

$d = $word.ActiveDocument

# --- locate the field whose code contains the m:'...'.fromHTMLURI() call ---
$targetField = $null
foreach ($fld in $d.Fields) {
    if ($fld.Code.Text -like "*fromHTMLURI*") {
        $targetField = $fld
    }
}
if ($targetField -eq $null) {
    throw "could not find the fromHTMLURI field"
}

# --- find the paragraph that contains the whole field (begin..end) ---
$fieldStart = $targetField.Code.Start - 1   # Code starts just after fldChar begin
$fieldEndCode = $targetField.Code.End       # Code ends just before the closing instrText/fldChar end

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Start -le $fieldStart -and $p.Range.End -gt $fieldEndCode) {
        $targetPara = $p
    }
}
if ($targetPara -eq $null) {
    throw "could not find the paragraph containing the field"
}

$paraRange = $targetPara.Range
# Range covering the whole paragraph's content, excluding the trailing
# paragraph mark, i.e. from the fldChar begin through the fldChar end.
$replaceRange = $d.Range($paraRange.Start, $paraRange.End - 1)

# Preserve the paragraph's own attributes (as found in the source document)
# so only its *content* changes, not the <w:p> element itself.
$paraAttrs = 'w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F"'

$apos = "'"

$runsXml = (
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>' + $apos + '</w:t></w:r>' +
    '<w:r><w:t>doc.html</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>' + $apos + '.fromHTMLURI()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>'
)

$xml = (
    '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
      '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
          '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:body>' +
              '<w:p ' + $paraAttrs + '>' + $runsXml + '</w:p>' +
            '</w:body>' +
          '</w:document>' +
        '</pkg:xmlData>' +
      '</pkg:part>' +
    '</pkg:package>'
)

$replaceRange.InsertXML($xml)
